# Stillingskoder.xlsx - add a calculated "førstestilling" column to Table1
#
# Adds a 3rd table column (C) that flags whether the Norwegian word
# "førstestilling" occurs in the stillingstittel (B) text, using:
#   =ISNUMBER(SEARCH("førstestilling",Table1[[#This Row],[stillingstittel]]))

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table (Table1) by one column - this extends ref/autoFilter from
# A1:B23 to A1:C23 and adds a new tableColumn entry.
$tbl = $ws.ListObjects.Item(1)
$col = $tbl.ListColumns.Add()

# Header text for the new column (also registers the shared string and
# renames the ListColumn to match, since the table column name tracks the
# header cell's text).
$ws.Range("C1").Value = "førstestilling"

# Fill in the calculated-column formula for every data row (2-23).
$formula = '=ISNUMBER(SEARCH("førstestilling",Table1[[#This Row],[stillingstittel]]))'
for ($r = 2; $r -le 23; $r++) {
    $ws.Range("C$r").Formula = $formula
}

# Match the new column's width to the rest of the table's look & feel.
$ws.Columns.Item(3).ColumnWidth = 16

# Reflect the author's post-edit selection/view state.
$ws.Range("C2").Select()

Write-Output "Added calculated column 'førstestilling' to Table1 (C2:C23)"
